$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1526.375
$ws.Range("I32").Value = 1653.6666
$ws.Range("J32").Value = 1450
$ws.Range("K32").Value = 1653.6666
$ws.Range("L32").Value = 1450
$ws.Range("M32").Value = -1327.6666
$ws.Range("N32").Value = -2102

# ARM row 35
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 30000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 30000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 30000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -30812

# ARM row 36
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 6735.154
$ws.Range("I36").Value = 4079.5715
$ws.Range("J36").Value = 9833.333000000001
$ws.Range("K36").Value = 4079.5715
$ws.Range("L36").Value = 9833.333000000001
$ws.Range("M36").Value = -3733.5715
$ws.Range("N36").Value = -10525.333

# ARM row 58
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 23000
$ws.Range("J58").Value = 23000
$ws.Range("L58").Value = 23000
$ws.Range("N58").Value = -23860

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 84518.836
$ws.Range("I74").Value = 126014
$ws.Range("J74").Value = 1528.5
$ws.Range("K74").Value = 126014
$ws.Range("L74").Value = 1528.5
$ws.Range("M74").Value = -125140
$ws.Range("N74").Value = -3276.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 84518.836
$ws.Range("I77").Value = 126014
$ws.Range("J77").Value = 1528.5
$ws.Range("K77").Value = 630070
$ws.Range("L77").Value = 7642.5
$ws.Range("M77").Value = -625702
$ws.Range("N77").Value = -16378.5

# ARM row 95
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 26000
$ws.Range("J95").Value = 26000
$ws.Range("L95").Value = 26000
$ws.Range("N95").Value = -31492

# BSM row 36
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 3110.7778
$ws.Range("I36").Value = 1249.625
$ws.Range("K36").Value = 1249.625
$ws.Range("M36").Value = -715.625

# BSM row 69
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 29800
$ws.Range("J69").Value = 29800
$ws.Range("L69").Value = 29800
$ws.Range("N69").Value = -31422

# BSM row 72
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H72").Value = 29800
$ws.Range("J72").Value = 29800
$ws.Range("L72").Value = 89400
$ws.Range("N72").Value = -97512

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 993.3393
$ws.Range("I31").Value = 819.8333
$ws.Range("J31").Value = 1123.4688
$ws.Range("K31").Value = 819.8333
$ws.Range("L31").Value = 1123.4688
$ws.Range("M31").Value = -524.8333
$ws.Range("N31").Value = -1713.4688

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 993.3393
$ws.Range("I34").Value = 819.8333
$ws.Range("J34").Value = 1123.4688
$ws.Range("K34").Value = 819.8333
$ws.Range("L34").Value = 1123.4688
$ws.Range("M34").Value = -617.8333
$ws.Range("N34").Value = -1527.4688

# CRP row 82
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 17727
$ws.Range("J82").Value = 24090.5
$ws.Range("L82").Value = 24090.5
$ws.Range("N82").Value = -24812.5

# CRP row 85
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H85").Value = 17727
$ws.Range("J85").Value = 24090.5
$ws.Range("L85").Value = 24090.5
$ws.Range("N85").Value = -26586.5

# CUL row 10
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 10363.3
$ws.Range("I10").Value = 14333.286
$ws.Range("J10").Value = 1100
$ws.Range("K10").Value = 42999.858
$ws.Range("L10").Value = 3300
$ws.Range("M10").Value = -42860.858
$ws.Range("N10").Value = -3578

# CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 179.70589
$ws.Range("I40").Value = 179.70589
$ws.Range("K40").Value = 718.82356
$ws.Range("M40").Value = -649.82356

# CUL row 62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4200
$ws.Range("J62").Value = 4333.3335
$ws.Range("L62").Value = 13000.0005
$ws.Range("N62").Value = -14372.0005

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 6716
$ws.Range("I63").Value = 3012
$ws.Range("K63").Value = 9036
$ws.Range("M63").Value = -8287

# CUL row 65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 4200
$ws.Range("J65").Value = 4333.3335
$ws.Range("L65").Value = 39000.0015
$ws.Range("N65").Value = -45864.0015

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 6716
$ws.Range("I66").Value = 3012
$ws.Range("K66").Value = 27108
$ws.Range("M66").Value = -23364

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 5123.3076
$ws.Range("I70").Value = 3350.5
$ws.Range("J70").Value = 6642.857
$ws.Range("K70").Value = 10051.5
$ws.Range("L70").Value = 19928.571
$ws.Range("M70").Value = -9736.5
$ws.Range("N70").Value = -20558.571

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 5123.3076
$ws.Range("I73").Value = 3350.5
$ws.Range("J73").Value = 6642.857
$ws.Range("K73").Value = 10051.5
$ws.Range("L73").Value = 19928.571
$ws.Range("M73").Value = -8959.5
$ws.Range("N73").Value = -22112.571

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3231
$ws.Range("I81").Value = 2013
$ws.Range("J81").Value = 3840
$ws.Range("K81").Value = 6039
$ws.Range("L81").Value = 11520
$ws.Range("M81").Value = -4916
$ws.Range("N81").Value = -13766

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 3231
$ws.Range("I84").Value = 2013
$ws.Range("J84").Value = 3840
$ws.Range("K84").Value = 18117
$ws.Range("L84").Value = 34560
$ws.Range("M84").Value = -12501
$ws.Range("N84").Value = -45792

# CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 11230510
$ws.Range("I114").Value = 7692630.5
$ws.Range("J114").Value = 22728620
$ws.Range("K114").Value = 23077891.5
$ws.Range("L114").Value = 68185860
$ws.Range("M114").Value = -23074637.5
$ws.Range("N114").Value = -68192368

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 383.66666
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 919.18555
$ws.Range("J131").Value = 919.18555
$ws.Range("L131").Value = 2757.55665
$ws.Range("N131").Value = -12837.55665

# GSM row 26
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 12500
$ws.Range("J26").Value = 12500
$ws.Range("L26").Value = 12500
$ws.Range("N26").Value = -13060

# GSM row 50
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 12500
$ws.Range("J50").Value = 12500
$ws.Range("L50").Value = 12500
$ws.Range("N50").Value = -13496

# GSM row 68
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4189
$ws.Range("N68").ClearContents()

# GSM row 71
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -10944
$ws.Range("N71").ClearContents()

# GSM row 117
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -31884

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2491.0557
$ws.Range("I132").Value = 2285.6843
$ws.Range("J132").Value = 2978.8125
$ws.Range("K132").Value = 6857.0529
$ws.Range("L132").Value = 8936.4375
$ws.Range("M132").Value = -4327.0529
$ws.Range("N132").Value = -13996.4375

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3160.7
$ws.Range("I122").Value = 3200
$ws.Range("J122").Value = 3134.5
$ws.Range("K122").Value = 9600
$ws.Range("L122").Value = 9403.5
$ws.Range("M122").Value = -7150
$ws.Range("N122").Value = -14303.5

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2989.5
$ws.Range("I132").Value = 2729.9736
$ws.Range("J132").Value = 4222.25
$ws.Range("K132").Value = 8189.9208
$ws.Range("L132").Value = 12666.75
$ws.Range("M132").Value = -5659.9208
$ws.Range("N132").Value = -17726.75

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1725
$ws.Range("I136").Value = 1685.7142
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5057.142599999999
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2507.142599999999
$ws.Range("N136").Value = -11100

# WVR row 57
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 24000.053
$ws.Range("J57").Value = 24000.053
$ws.Range("L57").Value = 24000.053
$ws.Range("N57").Value = -25508.053

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2152.3157
$ws.Range("I136").Value = 1589.9
$ws.Range("J136").Value = 2777.2222
$ws.Range("K136").Value = 4769.700000000001
$ws.Range("L136").Value = 8331.6666
$ws.Range("M136").Value = -2219.700000000001
$ws.Range("N136").Value = -13431.6666
